$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C10 ("Integer min" column, rule R30 row) changes from 18 to 1
$ws.Range("C10").Value = 1
